$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Title (appears twice: H1 heading + bold run near the end) - identical replacement both times
Replace-Text "Play Hellcatraz for Free - Review of Retro-Styled Slot Game" "Play Hellcatraz Slot for Free"

# "What we like" bullets
Replace-Text "High amount of paylines for increased winning potential" "High number of paylines for increased winning potential"
Replace-Text "Exciting bonus features including Super Free Spin and Lockdown Bonus" "Unique 8-bit graphics offer a retro arcade-style experience"
Replace-Text "Retro-style graphics offer a unique arcade-style experience" "Exciting bonus features including Super Free Spin and Cascading Reels"
Replace-Text "Wide range of betting options for all types of players" "Lockdown Bonus feature triggers free spins for big wins"

# "What we don't like" bullets
Replace-Text "Graphics may not appeal to all players" "Graphics may not appeal to new players"
Replace-Text "May seem overwhelming to new players due to high number of paylines" "Limited betting range from 0.20 to 20 coins"

# Meta title / description (italic) at the very end
Replace-Text "Read our review of Hellcatraz, a retro-style slot game with exciting bonus features. Play for free and maximize your winnings with Super Free Spin and Lockdown Bonus." "Get ready to play Hellcatraz, a retro-inspired slot game with exciting bonus features, for free."
